$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'309.67"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'1.17%"
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'35.59"
$ws.Range("D3").Style = "Normal"
$ws.Range("D4").Value = "'5.109"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'1.32%"
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'0.08195"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'3.76%"
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'2.047"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'-9.67%"
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'7.971"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'-0.34%"
$ws.Range("E7").Style = "Normal"
$ws.Range("B8").Value = "'GateToken"
$ws.Range("B8").Style = "Normal"
$ws.Range("C8").Value = "'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("C8").Style = "Normal"
$ws.Range("D8").Value = "'4.128"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'-0.62%"
$ws.Range("E8").Style = "Normal"
$ws.Range("B9").Value = "'BTSEToken"
$ws.Range("B9").Style = "Normal"
$ws.Range("C9").Value = "'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("C9").Style = "Normal"
$ws.Range("D9").Value = "'2.896"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'8.86%"
$ws.Range("E9").Style = "Normal"
$ws.Range("B10").Value = "'MXToken"
$ws.Range("B10").Style = "Normal"
$ws.Range("C10").Value = "'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("C10").Style = "Normal"
$ws.Range("D10").Value = "'0.9289"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'0.11%"
$ws.Range("E10").Style = "Normal"
$ws.Range("B11").Value = "'LiechtensteinCryptoassetsExchange"
$ws.Range("B11").Style = "Normal"
$ws.Range("C11").Value = "'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("C11").Style = "Normal"
$ws.Range("D11").Value = "'0.1084"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'10.64%"
$ws.Range("E11").Style = "Normal"
$ws.Range("B12").Value = "'WazirX"
$ws.Range("B12").Style = "Normal"
$ws.Range("C12").Value = "'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("C12").Style = "Normal"
$ws.Range("D12").Value = "'0.1922"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'2.84%"
$ws.Range("E12").Style = "Normal"
$ws.Range("B13").Value = "'MandalaExchangeToken"
$ws.Range("B13").Style = "Normal"
$ws.Range("C13").Value = "'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("C13").Style = "Normal"
$ws.Range("D13").Value = "'0.09510"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'5.62%"
$ws.Range("E13").Style = "Normal"
$ws.Range("B14").Value = "'BitrueCoin"
$ws.Range("B14").Style = "Normal"
$ws.Range("C14").Value = "'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("C14").Style = "Normal"
$ws.Range("D14").Value = "'0.03583"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'-4.47%"
$ws.Range("E14").Style = "Normal"
$ws.Range("B15").Value = "'BitMartToken"
$ws.Range("B15").Style = "Normal"
$ws.Range("C15").Value = "'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("C15").Style = "Normal"
$ws.Range("D15").Value = "'0.09903"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'-0.08%"
$ws.Range("E15").Style = "Normal"
$ws.Range("B16").Value = "'BitForexToken"
$ws.Range("B16").Style = "Normal"
$ws.Range("C16").Value = "'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("C16").Style = "Normal"
$ws.Range("D16").Value = "'0.001428"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'-0.12%"
$ws.Range("E16").Style = "Normal"
$ws.Range("B17").Value = "'TigerCash"
$ws.Range("B17").Style = "Normal"
$ws.Range("C17").Value = "'https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("C17").Style = "Normal"
$ws.Range("D17").Value = "'0.005727"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'0.13%"
$ws.Range("E17").Style = "Normal"
$ws.Range("B18").Value = "'LEO"
$ws.Range("B18").Style = "Normal"
$ws.Range("C18").Value = "'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("C18").Style = "Normal"
$ws.Range("D18").Value = "'3.475"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'0.36%"
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'0.3456"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'2.63%"
$ws.Range("E19").Style = "Normal"
$ws.Range("E20").Value = "'-0.71%"
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'5.102"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'0.43%"
$ws.Range("E21").Style = "Normal"
$ws.Range("D23").Value = "'0.04558"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'-0.41%"
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'0.001227"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'-0.57%"
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'0.004794"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'0.39%"
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'0.0001252"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'-3.75%"
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = "'0.0004451"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'-6.07%"
$ws.Range("E27").Style = "Normal"
$ws.Range("E39").Value = "'2.64%"
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'0.04900"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'-0.57%"
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'0.007668"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'-2.00%"
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'0.009853"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'26.06%"
$ws.Range("E42").Style = "Normal"
$ws.Range("E43").Value = "'-0.45%"
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'0.002118"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'-3.36%"
$ws.Range("E44").Style = "Normal"
$ws.Range("E45").Value = "'1.31%"
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'0.00006524"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'6.24%"
$ws.Range("E46").Style = "Normal"
$ws.Range("E47").Value = "'-0.01%"
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'175.39"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'238.84%"
$ws.Range("E48").Style = "Normal"
$ws.Range("E49").Value = "'-16.83%"
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "'0.00002101"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'-0.01%"
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = "'0.0002001"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'-0.01%"
$ws.Range("E51").Style = "Normal"
